$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '91.197.88'
$ws.Cells.Item(2,5).Value = '  +3.39%  '

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '3.150.41'
$ws.Cells.Item(3,5).Value = '  +2.32%  '

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '0.999'
$ws.Cells.Item(4,5).Value = '  -0.09%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '220.86'
$ws.Cells.Item(5,5).Value = '  +5.99%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '624.65'
$ws.Cells.Item(6,5).Value = '  +0.43%  '

$ws.Cells.Item(7,5).Value = '  +2.25%  '

$ws.Cells.Item(8,5).Value = '  +10.46%  '

$ws.Cells.Item(9,5).Value = '  +0.01%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '3.146.31'
$ws.Cells.Item(10,5).Value = '  +2.25%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.730'
$ws.Cells.Item(11,5).Value = '  +21.93%  '

$ws.Cells.Item(12,5).Value = '  +6.21%  '

$ws.Cells.Item(13,5).Value = '  +7.24%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '34.51'
$ws.Cells.Item(14,5).Value = '  +9.27%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '5.43'
$ws.Cells.Item(15,5).Value = '  +3.20%  '

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '91.283.68'
$ws.Cells.Item(16,5).Value = '  +3.90%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '3.720.63'
$ws.Cells.Item(17,5).Value = '  +1.99%  '

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '3.139.92'
$ws.Cells.Item(18,5).Value = '  +1.99%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '3.86'
$ws.Cells.Item(19,5).Value = '  +21.14%  '

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '0.0000229'
$ws.Cells.Item(20,5).Value = '  +8.59%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '14.28'
$ws.Cells.Item(21,5).Value = '  +8.86%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '434.13'
$ws.Cells.Item(22,5).Value = '  +3.19%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '8.79'
$ws.Cells.Item(23,5).Value = '  +7.98%  '

$ws.Cells.Item(24,5).Value = '  +6.74%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '6.09'
$ws.Cells.Item(25,5).Value = '  +11.91%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '12.50'
$ws.Cells.Item(26,5).Value = '  +8.90%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '83.69'
$ws.Cells.Item(27,5).Value = '  +2.39%  '

$ws.Cells.Item(28,5).Value = '  -0.25%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '0.168'
$ws.Cells.Item(29,5).Value = '  +6.66%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '9.02'
$ws.Cells.Item(30,5).Value = '  +12.45%  '

$ws.Cells.Item(31,2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(31,3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '0.997'
$ws.Cells.Item(31,5).Value = '  -8.34%  '

$ws.Cells.Item(32,2).Value = 'Bittensor'
$ws.Cells.Item(32,3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '535.34'
$ws.Cells.Item(32,5).Value = '  +6.18%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '3.95'
$ws.Cells.Item(33,5).Value = '  +11.97%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '7.33'
$ws.Cells.Item(34,5).Value = '  +11.49%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.33'
$ws.Cells.Item(35,5).Value = '  +8.41%  '

$ws.Cells.Item(36,5).Value = '  +6.03%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '23.61'
$ws.Cells.Item(37,5).Value = '  +6.22%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '1.87'
$ws.Cells.Item(38,5).Value = '  +3.39%  '

$ws.Cells.Item(39,5).Value = '  +0.44%  '

$ws.Cells.Item(40,5).Value = '  +0.03%  '

$ws.Cells.Item(41,2).Value = 'Hedera'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.0790'
$ws.Cells.Item(41,5).Value = '  +18.34%  '

$ws.Cells.Item(42,2).Value = 'Stellar'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.146'
$ws.Cells.Item(42,5).Value = '  +7.03%  '

$ws.Cells.Item(43,2).Value = 'USDe'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '1.00'
$ws.Cells.Item(43,5).Value = '  +0.01%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '0.379'
$ws.Cells.Item(44,5).Value = '  +5.90%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '1.93'
$ws.Cells.Item(45,5).Value = '  +6.67%  '

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '143.86'
$ws.Cells.Item(46,5).Value = '  -3.33%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '44.20'
$ws.Cells.Item(47,5).Value = '  +1.88%  '

$ws.Cells.Item(48,5).Value = '  +11.28%  '

$ws.Cells.Item(49,2).Value = 'Aave'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '169.37'
$ws.Cells.Item(49,5).Value = '  +8.69%  '

$ws.Cells.Item(50,2).Value = 'FLOKI'
$ws.Cells.Item(50,3).Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.000265'
$ws.Cells.Item(50,5).Value = '  +24.17%  '

$ws.Cells.Item(51,2).Value = 'Mantle'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.746'
$ws.Cells.Item(51,5).Value = '  +6.56%  '

